$d = $word.ActiveDocument

# 1) "Code Refactorings done in the project:" - merge the three runs
#    (and drop the spell-check proofErr markers around "Refactorings")
#    into a single run by replacing the text with itself.
$d.Content.Find.Execute(
    "Code Refactorings done in the project:", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Code Refactorings done in the project:", 2) | Out-Null

# 2) "Edit Book: ... has the ability to ..." - merge the three runs
#    (and drop the grammar-check proofErr markers) into a single run.
$d.Content.Find.Execute(
    "Edit Book: The application now has the ability to edit the details of a book and update the same in the database",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Edit Book: The application now has the ability to edit the details of a book and update the same in the database", 2) | Out-Null

# 3) "There are various other details ..." - merge the three runs
#    (and drop the grammar-check proofErr markers) into a single run.
$apos = [char]0x2019
$dash = [char]0x2013
$closing = "There are various other details that I could have added, but couldn" + $apos + "t like " + $dash + " specification pattern especially in the filter criteria logic."
$d.Content.Find.Execute($closing, $true, $false, $false, $false, $false, $true, 1, $false, $closing, 2) | Out-Null

# Now append a blank paragraph plus a brand-new closing paragraph after it.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter() | Out-Null

$blank = $d.Paragraphs.Last
$blank.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "I have no experience in Vue.js, but I have tried my best to accommodate all the requirements, though styling has not been done to a great extent."
